$d = $word.ActiveDocument

# Locate the original sentence that is being retyped.
$rng = $d.Content
$found = $rng.Find.Execute("Những điểm có khả năng cạnh tranh thua kém, là ở các yếu tố")
if (-not $found) {
    throw "could not find target sentence"
}
$matchStart = $rng.Start

# Replace it with the corrected wording (comma after "kém" dropped, replaced
# by a simple space before "là ở các yếu tố").
$newText = "Những điểm có khả năng cạnh tranh thua kém là ở các yếu tố"
$rng.Text = $newText

# The trailing "_GoBack" bookmark (which used to sit at the very end of the
# paragraph) now has to sit right after "... thua kém " -- i.e. right where
# the author's cursor was left after retyping that chunk.
if ($d.Bookmarks.Exists("_GoBack")) {
    $bmOld = $d.Bookmarks("_GoBack")
    $bmOld.Delete()
}

$firstLen = "Những điểm có khả năng cạ".Length
$secondLen = "nh tranh thua kém ".Length
$thirdLen = "là ở các yếu tố".Length

$bmPos = $matchStart + $firstLen + $secondLen
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Force the paragraph's run list to split at the same three boundaries the
# author's retyping produced, toggling a formatting property on and back off
# so the final character formatting is unchanged but a run break is created.
$split1 = $matchStart + $firstLen
$probe1 = $d.Range($matchStart, $split1)
$probe1.Font.Bold = 1
$probe1.Font.Bold = 0

$split2 = $matchStart + $firstLen + $secondLen + $thirdLen
$probe2 = $d.Range($matchStart, $split2)
$probe2.Font.Bold = 1
$probe2.Font.Bold = 0

Write-Output "done"
